$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new data rows right after the header block (before the former row 9),
# shifting all subsequent rows (old 9..89) down to (11..91).
$ws.Rows("9:10").Insert()

# --- New row 9 ---
$ws.Range("A9").Value = 5
$ws.Range("B9").Value = "Macroferia Regional de Talca"
$ws.Range("C9").Value = "Maule"
$ws.Range("D9").Value = 44613
$ws.Range("E9").Value = 7
$ws.Range("F9").Value = "Fruta"
$ws.Range("G9").Value = 100103
$ws.Range("H9").Value = "Frutos de hueso (carozo)"
$ws.Range("I9").Value = 100103002
$ws.Range("J9").Value = "Ciruela"
$ws.Range("K9").Value = "Black Amber"
$ws.Range("L9").Value = "Especial"
$ws.Range("M9").Value = 200
$ws.Range("N9").Value = 12000
$ws.Range("O9").Value = 12000
$ws.Range("P9").Value = 12000
$ws.Range("Q9").Value = "`$/bandeja 18 kilos granel"
$ws.Range("R9").Value = "Región de O'Higgins"
$ws.Range("S9").Value = 667
$ws.Range("T9").Value = 18

# --- New row 10 ---
$ws.Range("A10").Value = 5
$ws.Range("B10").Value = "Macroferia Regional de Talca"
$ws.Range("C10").Value = "Maule"
$ws.Range("D10").Value = 44613
$ws.Range("E10").Value = 7
$ws.Range("F10").Value = "Fruta"
$ws.Range("G10").Value = 100103
$ws.Range("H10").Value = "Frutos de hueso (carozo)"
$ws.Range("I10").Value = 100103002
$ws.Range("J10").Value = "Ciruela"
$ws.Range("K10").Value = "Black Amber"
$ws.Range("L10").Value = "Primera"
$ws.Range("M10").Value = 200
$ws.Range("N10").Value = 10000
$ws.Range("O10").Value = 10000
$ws.Range("P10").Value = 10000
$ws.Range("Q10").Value = "`$/bandeja 18 kilos granel"
$ws.Range("R10").Value = "Región de O'Higgins"
$ws.Range("S10").Value = 556
$ws.Range("T10").Value = 18

$dim = $ws.UsedRange.Address()
Write-Host "Final UsedRange: $dim"
